$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E3").Value = "['Normal', 'ParamViolation']"

# Row 38
$ws.Range("D38").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['Normal', 'HardwareFault']"

# Row 39
$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"

# Row 56
$ws.Range("D56").Value = "[0, 0, 0, 0, 0, 0, 0]"
$ws.Range("E56").Value = "[]"

# Row 73
$ws.Range("D73").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'HardwareFault']"

# Row 83
$ws.Range("D83").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E83").Value = "['Normal', 'SurroundingEnvironment']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 0, 0, 0, 0, 1]"
$ws.Range("E113").Value = "['Normal', 'SoftwareFault']"
